# Add "Wins" / "Losses" / "Ties" season-record columns (AD:AF) to the
# DET_2010 sheet, matching the header style already used by the other
# header cells (bold font, thin border, centered/top alignment) and
# filling every data row (2-43) with the team's 81-81-0 record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the exact header style (bold/border/centered) from an existing
# header cell instead of re-building it by hand.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-43): 81 wins, 81 losses, 0 ties -------------------
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD
    $ws.Cells.Item($r, 31).Value = 81  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}

Write-Output "season record columns added"
